$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assign "Jabesi" as assignee for the tasks the author picked up,
# and record the estimated time (in hours) to finish each.
$ws.Range("D4").Value = "Jabesi"
$ws.Range("E4").Value = 14

$ws.Range("D6").Value = "Jabesi"
$ws.Range("E6").Value = 4

$ws.Range("D7").Value = "Jabesi"
$ws.Range("E7").Value = 4

# Move the active selection to reflect where the author left off.
[void]$ws.Range("E5").Select()
